$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style of L4 ("Unique") onto the two new header cells
# so M4/N4 ("Pattern", "Pattern Type") look like the rest of the header row.
$ws.Range("L4").Copy()
$ws.Range("M4:N4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("M4").Value = "Pattern"
$ws.Range("N4").Value = "Pattern Type"

$ws.Range("M4:N4").Select()
